$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("H2").Value = 0.5797639999999999
$ws.Range("M2").Value = 52.36202133333333
$ws.Range("N2").Value = 157.086064
$ws.Range("O2").Value = 0.5039383032147022
$ws.Range("P2").Value = 0.5039383032147023
$ws.Range("Q2").Value = 10.11920497876622
$ws.Range("R2").Value = 91.07284480889599
$ws.Range("S2").Value = 0.5039383032147022
$ws.Range("T2").Value = 0.5039383032147023

# Row 3
$ws.Range("H3").Value = 0.5797639999999999
$ws.Range("O3").Value = 0.1251951767241128
$ws.Range("P3").Value = 0.1251951767241128
$ws.Range("Q3").Value = 2.513949917167555
$ws.Range("S3").Value = 0.1251951767241128
$ws.Range("T3").Value = 0.1251951767241128

# Row 4
$ws.Range("H4").Value = 0.5797639999999999
$ws.Range("M4").Value = 10.240131
$ws.Range("N4").Value = 30.720393
$ws.Range("O4").Value = 0.09855223517796473
$ws.Range("P4").Value = 0.09855223517796474
$ws.Range("Q4").Value = 1.978953103028
$ws.Range("R4").Value = 17.810577927252
$ws.Range("S4").Value = 0.09855223517796473
$ws.Range("T4").Value = 0.09855223517796474

# Row 5
$ws.Range("H5").Value = 0.5797639999999999
$ws.Range("M5").Value = 1.816419
$ws.Range("N5").Value = 5.449257
$ws.Range("O5").Value = 0.01748143187520975
$ws.Range("P5").Value = 0.01748143187520975
$ws.Range("Q5").Value = 0.351031448372
$ws.Range("R5").Value = 3.159283035348
$ws.Range("S5").Value = 0.01748143187520975
$ws.Range("T5").Value = 0.01748143187520975

# Row 6
$ws.Range("H6").Value = 0.5797639999999999
$ws.Range("M6").Value = 15.88623066666667
$ws.Range("N6").Value = 47.658692
$ws.Range("O6").Value = 0.1528909679722582
$ws.Range("P6").Value = 0.1528909679722583
$ws.Range("Q6").Value = 3.070088212076445
$ws.Range("R6").Value = 27.630793908688
$ws.Range("S6").Value = 0.1528909679722582
$ws.Range("T6").Value = 0.1528909679722583

# Row 7
$ws.Range("H7").Value = 0.5797639999999999
$ws.Range("M7").Value = 10.59233466666667
$ws.Range("N7").Value = 31.777004
$ws.Range("O7").Value = 0.1019418850357522
$ws.Range("P7").Value = 0.1019418850357522
$ws.Range("Q7").Value = 2.047018105228445
$ws.Range("R7").Value = 18.423162947056
$ws.Range("S7").Value = 0.1019418850357522
$ws.Range("T7").Value = 0.1019418850357522
